# Apply the Backlog.xlsx edit:
#  - On the ARCHIVE sheet, update status (column E) for rows 25-30
#    from "In Progress" to "Complete".
#  - Update the sheet's active selection to E24:E30 (active cell E24),
#    which also drops the stale topLeftCell="A15" scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHIVE")
$ws.Activate()

$ws.Range("E25:E30").Value = "Complete"

$ws.Range("E24:E30").Select()
